# Apply trade #31 closing updates to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.77
$summary.Range("B4").Value = -0.23
$summary.Range("B5").Value = -0.15
$summary.Range("B6").Value = 31
$summary.Range("B8").Value = 16
$summary.Range("B9").Value = 29.03

# --- Sheet: Strategy Status (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.77
$status.Range("D4").Value = 31
$status.Range("E4").Value = -0.23
$status.Range("F4").Value = -0.23
$status.Range("G4").Value = 29.03

# --- New trade row shared by "All Trades" and "MarketMaking" sheets ---
$tradeNum = 31
$date = "2026-02-17"
$time = "15:22:31"
$strategy = "MarketMaking"
$side = "UP"
$entryPrice = 0.75
$exitPrice = 0.68
$status2 = "CLOSED"
$pnlPct = -9.333299999999999
$pnlDollar = -0.07000000000000001
$capitalAfter = 99.77
$entrySlip = 0
$exitSlip = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.14

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 32
    $ws.Cells.Item($row, 1).Value = $tradeNum
    # Force the date column to stay as literal text (matches existing rows,
    # which store the date as a plain string, not an Excel date serial).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $date
    $ws.Cells.Item($row, 3).Value = $time
    $ws.Cells.Item($row, 4).Value = $strategy
    $ws.Cells.Item($row, 5).Value = $side
    $ws.Cells.Item($row, 6).Value = $entryPrice
    $ws.Cells.Item($row, 7).Value = $exitPrice
    $ws.Cells.Item($row, 8).Value = $status2
    $ws.Cells.Item($row, 9).Value = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlip
    $ws.Cells.Item($row, 13).Value = $exitSlip
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = $exitReason
    $ws.Cells.Item($row, 17).Value = $duration
}

$wb.Save()
